$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON) meanEMG values for columns B:E
$ws.Range("B2").Value = 282.48068171035595
$ws.Range("C2").Value = 262.95206511980251
$ws.Range("D2").Value = 282.36404105544182
$ws.Range("E2").Value = 259.9637642171167

# Update row 3 (STR) meanEMG values for columns B:E
$ws.Range("B3").Value = 307.03862696799189
$ws.Range("C3").Value = 259.5601453432763
$ws.Range("D3").Value = 315.56734272789538
$ws.Range("E3").Value = 259.22057767643787

# Update the selection to match the new selected range B1:E3
$ws.Range("B1:E3").Select()
